# Add a new "test_name" / "Test Name" row to the column-name mapping sheet
# (mirrors the addition of a new row to the manuscript's summary column
# name lookup table).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "test_name"
$ws.Range("B8").Value = "Test Name"

# Match the author's final cursor position recorded in the saved file.
$ws.Range("C13").Select() | Out-Null
